$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Members")

# New headers for the "Father"/"Mother" cross-reference columns
$ws.Range("E1").Value = "Father"
$ws.Range("F1").Value = "Mother"
$ws.Range("E1:F1").Font.Bold = $true

# Lea's parents: Paul (father) and Isa (mother)
$ws.Range("E5").Value = "Paul"
$ws.Range("F5").Value = "Isa"

# Elias' parents: Paul (father) and Isa (mother)
$ws.Range("E6").Value = "Paul"
$ws.Range("F6").Value = "Isa"

# Dave, Clara, Bryan and Fiona's father: Elias
$ws.Range("E8").Value = "Elias"
$ws.Range("E9").Value = "Elias"
$ws.Range("E10").Value = "Elias"
$ws.Range("E11").Value = "Elias"

# Alain's parents: Dave (father) and Katell (mother)
$ws.Range("E12").Value = "Dave"
$ws.Range("F12").Value = "Katell"

$ws.Range("E13").Select()
